$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value of D1 - the color is now green (00ff00) instead of black (000000)
$ws.Range("D1").Value = "coloR=   00ff00"

# Move the selection to E7 (last place the user clicked)
$ws.Range("E7").Select()
